$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.366703152656555
$ws.Range("B1").Value = 3.826034784317017
$ws.Range("C1").Value = 3.544075727462769
$ws.Range("D1").Value = 3.3410804271698
$ws.Range("E1").Value = 1.219692707061768
